$d = $word.ActiveDocument

function Merge-Runs($firstRunText, $fullText) {
    # Locate the first run's text to find the cut point (end of the run
    # that should keep its original rPr/formatting).
    $rngHead = $d.Range(0, $d.Content.End)
    $rngHead.Find.Execute($firstRunText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    $cut = $rngHead.End

    # Locate the whole (already-existing, just split across runs) text to
    # find where the paragraph's sentence ends.
    $rngFull = $d.Range(0, $d.Content.End)
    $rngFull.Find.Execute($fullText, $true, $false, $false, $false, $false, `
                           $true, 1, $false, "", 0)
    $fullEnd = $rngFull.End

    $tail = $fullText.Substring($firstRunText.Length)

    # Remove the remaining runs' text, then re-append it onto the first
    # run so the paragraph ends up with a single merged run.
    $rngDel = $d.Range($cut, $fullEnd)
    $rngDel.Delete()

    $insPoint = $d.Range($cut, $cut)
    $insPoint.InsertAfter($tail)
}

# 1) Merge the split "wop" runs in the precondition sentence.
Merge-Runs "Mindst en Workoutplan w" "Mindst en Workoutplan wop eksisterer."

# 2) Insert the new precondition paragraph right after it.
$precondPara = $d.Paragraphs.Item(11)
$precondPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(12)
$newPara.Range.Text = "P har mindst en wop tilknyttet."

# 3) Merge the split "wop" runs in the postcondition sentences.
Merge-Runs "En p valgte en w" "En p valgte en wop."
Merge-Runs "En w" "En wop blev præsenteret for bruger."
